# Update countries & provincias Spain
# Applies the "Pais" COVID-19 stats refresh:
#   - bumps the "Datos actualizados ..." timestamp in A1
#   - refreshes Casos totales / Nuevos casos / Casos activos / Recuperados /
#     Casos criticos / Muertes for the countries whose daily figures moved
#   - because some countries changed rank order with their new figures, a
#     handful of rows also need their country name (column A) swapped with a
#     neighbouring row to keep the table sorted by total cases

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Timestamp banner
$ws.Range("A1").Value = "Datos actualizados a 25 de Agosto de 2020 a las 17:42"

# Country name swaps caused by re-ranking (column A only)
$nameUpdates = @(
    @{ Row=45;  Name="Emiratos Arabes Unidos" },
    @{ Row=46;  Name="Paises Bajos" },
    @{ Row=94;  Name="Grecia" },
    @{ Row=95;  Name="Guayana Francesa" },
    @{ Row=143; Name="Jordania" },
    @{ Row=144; Name="Malta" },
    @{ Row=145; Name="Nueva Zelanda" },
    @{ Row=154; Name="Liberia" },
    @{ Row=155; Name="Togo" },
    @{ Row=157; Name="Trinidad yTobago" },
    @{ Row=158; Name="Niger" },
    @{ Row=174; Name="Papua Nueva Guinea" },
    @{ Row=175; Name="Comoras" },
    @{ Row=176; Name="Islas Feroe" },
    @{ Row=177; Name="San Martin (Parte Holandesa)" },
    @{ Row=214; Name="Islas Malvinas" },
    @{ Row=215; Name="Montserrat" }
)

foreach ($u in $nameUpdates) {
    $ws.Cells.Item($u.Row, 1).Value = $u.Name
}

# Updated statistics: B=Casos totales, C=Nuevos casos, D=Casos activos,
# E=Recuperados, G=Muertes hoy, H=Muertes  (F=Casos criticos stays 0 throughout)
$dataUpdates = @(
    @{ Row=4;   B=5922310; C=6680;  D=3219333; E=2521535; G=328; H=181442 },
    @{ Row=5;   B=3627961; C=744;   D=2778709; E=733776;  G=25;  H=115476 },
    @{ Row=6;   B=3193917; C=29036; D=2427097; E=707996;  G=278; H=58824 },
    @{ Row=13;  B=400985;  C=1417;  D=374463;  E=15564;   G=42;  H=10958 },
    @{ Row=15;  B=350867;  C=0;     D=263202;  E=80263;   G=36;  H=7402 },
    @{ Row=23;  B=236516;  C=399;   D=209600;  E=17579;   G=1;   H=9337 },
    @{ Row=27;  B=125747;  C=100;   D=111694;  E=4969;    G=1;   H=9084 },
    @{ Row=45;  B=67621;   C=339;   D=58754;   E=8490;    G=1;   H=377 },
    @{ Row=46;  B=67543;   C=415;   D=0;       E=0;       G=5;   H=6207 },
    @{ Row=68;  B=32803;   C=246;   D=19055;   E=13189;   G=5;   H=559 },
    @{ Row=90;  B=10426;   C=31;    D=9150;    E=1012;    G=0;   H=264 },
    @{ Row=94;  B=8987;    C=168;   D=3804;    E=4940;    G=1;   H=243 },
    @{ Row=95;  B=8875;    C=0;     D=8363;    E=456;     G=0;   H=56 },
    @{ Row=96;  B=8759;    C=154;   D=4530;    E=3970;    G=5;   H=259 },
    @{ Row=126; B=2971;    C=12;    D=2816;    E=143;     G=0;   H=12 },
    @{ Row=143; B=1716;    C=77;    D=1344;    E=358;     G=0;   H=14 },
    @{ Row=144; B=1705;    C=38;    D=1029;    E=666;     G=0;   H=10 },
    @{ Row=145; B=1690;    C=7;     D=1539;    E=129;     G=0;   H=22 },
    @{ Row=154; B=1295;    C=5;     D=821;     E=392;     G=0;   H=82 },
    @{ Row=155; B=1295;    C=0;     D=914;     E=354;     G=0;   H=27 },
    @{ Row=157; B=1184;    C=85;    D=165;     E=1004;    G=0;   H=15 },
    @{ Row=158; B=1172;    C=0;     D=1084;    E=19;      G=0;   H=69 },
    @{ Row=174; B=419;     C=18;    D=232;     E=183;     G=0;   H=4 },
    @{ Row=175; B=417;     C=0;     D=396;     E=14;      G=0;   H=7 },
    @{ Row=176; B=411;     C=1;     D=344;     E=67;      G=0;   H=0 },
    @{ Row=177; B=408;     C=0;     D=147;     E=244;     G=0;   H=17 },
    @{ Row=214; B=13;      C=0;     D=13;      E=0;       G=0;   H=0 },
    @{ Row=215; B=13;      C=0;     D=12;      E=0;       G=0;   H=1 }
)

foreach ($u in $dataUpdates) {
    $r = $u.Row
    $ws.Cells.Item($r, 2).Value = $u.B
    $ws.Cells.Item($r, 3).Value = $u.C
    $ws.Cells.Item($r, 4).Value = $u.D
    $ws.Cells.Item($r, 5).Value = $u.E
    $ws.Cells.Item($r, 7).Value = $u.G
    $ws.Cells.Item($r, 8).Value = $u.H
}
